$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 85: correct the date/time stamp in column A ---
$ws.Cells.Item(85, 1).Value2 = 45460.2916666667

# --- Row 86: new data row appended by the R script ---
# Column A (date) needs the same style (numFmt) as the rest of the date
# column, so copy it from the cell directly above rather than building a
# brand-new style entry.
$ws.Cells.Item(85, 1).Copy()
$ws.Cells.Item(86, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(86, 1).Value2 = 45461.6438310185

$ws.Cells.Item(86, 2).Value2 = 2700
$ws.Cells.Item(86, 3).Value2 = 6.1399998664856
$ws.Cells.Item(86, 4).Value2 = 6.07999992370605
$ws.Cells.Item(86, 5).Value2 = 6.11999988555908
$ws.Cells.Item(86, 6).Value2 = 6.07999992370605

# adj_close was written as text (matches the other rows, which all store
# this column's numeric-looking values as shared strings). Force text
# storage with a leading quote, then drop back to the default "Normal"
# style so no stray quote-prefixed style lingers on the cell.
$ws.Cells.Item(86, 7).Value = "'6.07999992370605"
$ws.Cells.Item(86, 7).Style = "Normal"

$ws.Cells.Item(86, 8).Value = "PAL.MI"
